$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: SubmittedDate -> ApprovedDate
$ws.Range("B1").Value = "ApprovedDate"

# Row 2: #00-10/fdgds/Food Expense/Debit Card -> #00-12/this is test/Travel Expense/Credit Card
$ws.Range("A2").Value = "#00-12"
$ws.Range("B2").Value = 43541.615219907406
$ws.Range("C2").Value = 43491.001388888886
$ws.Range("D2").Value = "this is test"
$ws.Range("F2").Value = "Travel Expense"
$ws.Range("G2").Value = "Credit Card"
$ws.Range("H2").Value = "Attached"
$ws.Range("I2").Value = 567.0

# Row 3: #00-11/tyw/Food Expense/Credit Card/Attached -> #00-16/effective/Food Expense/Credit Card/Not Attached
$ws.Range("A3").Value = "#00-16"
$ws.Range("B3").Value = 43541.63195601852
$ws.Range("C3").Value = 43493.001388888886
$ws.Range("D3").Value = "effective"
$ws.Range("F3").Value = "Food Expense"
$ws.Range("G3").Value = "Credit Card"
$ws.Range("H3").Value = "Not Attached"
$ws.Range("I3").Value = 6576.0

# Row 4: #00-12/this is test/Travel Expense/Credit Card -> #00-17/effective two/Food Expense/Credit Card
$ws.Range("A4").Value = "#00-17"
$ws.Range("B4").Value = 43541.61491898148
$ws.Range("C4").Value = 43472.00208333333
$ws.Range("D4").Value = "effective two"
$ws.Range("F4").Value = "Food Expense"
$ws.Range("G4").Value = "Credit Card"
$ws.Range("H4").Value = "Attached"
$ws.Range("I4").Value = 6576.0

# Row 5 (#00-16/effective/...) is removed entirely.
$ws.Range("A5:I5").Value = ""

# Row 6 (#00-17/effective two/...) becomes the Total row: every column blank
# except Expense Type = "Total" and Amount = sum of the remaining amounts.
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = "Total"
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = 13719.0
